# Burndown template update:
#  - adjust Sprint 3 (E column) actuals/remaining on "Burndown Chart"
#  - tweak the Remaining-effort formulas so the burndown line reflects the new totals
#  - on "User Stories": mark US3 as done, add two "optional-points" helper formulas,
#    flag US6 as "Optional", extend Table1 to include it, and add a
#    "*total not including optional" footer row
#  - view-state: update the active-cell selections to match where the authors left off

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Burndown Chart"
$ws2 = $wb.Worksheets.Item(2)   # "User Stories"

# ---------------------------------------------------------------------------
# 1) User Stories sheet: data edits first (so Burndown Chart totals below
#    are already consistent with the new Points total)
# ---------------------------------------------------------------------------

# US3 (row 4) is now marked Done
$ws2.Range("D4").Value = "y"

# Remove the old helper formula in G4 (superseded by the F-column helpers below)
$ws2.Range("G4").ClearContents()

# New helper column F: running "points so far" style computations
$ws2.Range("F2").Value = 3
$ws2.Range("F3").Formula = "=SUM(Table1[[#This Row],[Points]],B5,B6)"
$ws2.Range("F4").Formula = "=Table1[[#This Row],[Points]]"
$ws2.Range("F5").Formula = "=SUM(F2:F4)"

# US6 (row 7) is optional work, not part of the committed total
$ws2.Range("C7").Value = "Optional"

# Total: now excludes the optional US6 row (B2:B6 instead of B2:B7)
$ws2.Range("B8").Formula = "=SUM(B2:B6)"

# Extend Table1 (US/Points/Column1/Done?) down to the new row 9
$tbl = $ws2.ListObjects.Item(1)
$tbl.Resize($ws2.Range("A1:D9"))

# New footnote row 9, matching the formatting of the row above it
$ws2.Range("A7").Copy()
$ws2.Range("A9").PasteSpecial(-4122)
$ws2.Range("D7").Copy()
$ws2.Range("D9").PasteSpecial(-4122)
$ws2.Range("C7").Copy()
$ws2.Range("C9").PasteSpecial(-4122)

$ws2.Range("A7").Copy()
$ws2.Range("B9").PasteSpecial(-4122)
$ws2.Range("B9").Value = "*total not including optional"
$ws2.Range("B9").Font.Color = 255

# ---------------------------------------------------------------------------
# 2) Burndown Chart sheet: Sprint 3 actuals/remaining change with the new totals
# ---------------------------------------------------------------------------

$ws1.Range("E2").Value = 18
$ws1.Range("E3").Value = 10

$ws1.Range("B5").Value = 31
$ws1.Range("C5").Formula = "=`$B`$5-SUM(C3)"
$ws1.Range("D5").Formula = "=`$B`$5-SUM(C3:D3)"
$ws1.Range("E5").Formula = "=`$B`$5-SUM(C3:E3)"

# ---------------------------------------------------------------------------
# 3) View state: leave the selections where the authors left them
# ---------------------------------------------------------------------------

$ws2.Range("B9").Select()
$ws1.Range("E15").Select()
